# Replace the "ser: 111" blog entry with a new "ser: 114" blog entry.
# (Sheet1!I7 held the "type: blog / width: 2 / height: 1 / ser: 111" widget
# config; the author swapped it for entry 114.)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I7").Value = "type: blog" + "`n" + "width: 2" + "`n" + "height: 1" + "`n" + "ser: 114"
